# Natmi following Dr Hou advice
# Add a third sending/target cluster ("ECs") to the Comp-Itgb3 LR-pairs table
# and refresh all derived specificity/expression values for the 3x3 cluster grid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Comp"
$ws.Range("C2").Value = "Itgb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4181183333333333
$ws.Range("H2").Value = 1.254355
$ws.Range("I2").Value = 0.01571144052599341
$ws.Range("J2").Value = 0.01571144052599341
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.970048
$ws.Range("N2").Value = 26.910144
$ws.Range("O2").Value = 0.487108783009476
$ws.Range("P2").Value = 0.4871087830094759
$ws.Range("Q2").Value = 3.75054151968
$ws.Range("R2").Value = 33.75487367712
$ws.Range("S2").Value = 0.00765318067394241
$ws.Range("T2").Value = 0.007653180673942408

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Comp"
$ws.Range("C3").Value = "Itgb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.4181183333333333
$ws.Range("H3").Value = 1.254355
$ws.Range("I3").Value = 0.01571144052599341
$ws.Range("J3").Value = 0.01571144052599341
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.012070666666666
$ws.Range("N3").Value = 27.036212
$ws.Range("O3").Value = 0.489390778604016
$ws.Range("P3").Value = 0.489390778604016
$ws.Range("Q3").Value = 3.768111967028889
$ws.Range("R3").Value = 33.91300770325999
$ws.Range("S3").Value = 0.007689034112006604
$ws.Range("T3").Value = 0.007689034112006603

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Comp"
$ws.Range("C4").Value = "Itgb3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4181183333333333
$ws.Range("H4").Value = 1.254355
$ws.Range("I4").Value = 0.01571144052599341
$ws.Range("J4").Value = 0.01571144052599341
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4327576666666667
$ws.Range("N4").Value = 1.298273
$ws.Range("O4").Value = 0.02350043838650813
$ws.Range("P4").Value = 0.02350043838650813
$ws.Range("Q4").Value = 0.1809439143238889
$ws.Range("R4").Value = 1.628495228915
$ws.Range("S4").Value = 0.0003692257400443949
$ws.Range("T4").Value = 0.0003692257400443948

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Comp"
$ws.Range("C5").Value = "Itgb3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 24.89087966666667
$ws.Range("H5").Value = 74.672639
$ws.Range("I5").Value = 0.9353131502385497
$ws.Range("J5").Value = 0.9353131502385496
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.970048
$ws.Range("N5").Value = 26.910144
$ws.Range("O5").Value = 0.487108783009476
$ws.Range("P5").Value = 0.4871087830094759
$ws.Range("Q5").Value = 223.272385372224
$ws.Range("R5").Value = 2009.451468350016
$ws.Range("S5").Value = 0.4555992503454591
$ws.Range("T5").Value = 0.455599250345459

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Comp"
$ws.Range("C6").Value = "Itgb3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 24.89087966666667
$ws.Range("H6").Value = 74.672639
$ws.Range("I6").Value = 0.9353131502385497
$ws.Range("J6").Value = 0.9353131502385496
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.012070666666666
$ws.Range("N6").Value = 27.036212
$ws.Range("O6").Value = 0.489390778604016
$ws.Range("P6").Value = 0.489390778604016
$ws.Range("Q6").Value = 224.3183665114964
$ws.Range("R6").Value = 2018.865298603468
$ws.Range("S6").Value = 0.4577336308338188
$ws.Range("T6").Value = 0.4577336308338187

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Comp"
$ws.Range("C7").Value = "Itgb3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 24.89087966666667
$ws.Range("H7").Value = 74.672639
$ws.Range("I7").Value = 0.9353131502385497
$ws.Range("J7").Value = 0.9353131502385496
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4327576666666667
$ws.Range("N7").Value = 1.298273
$ws.Range("O7").Value = 0.02350043838650813
$ws.Range("P7").Value = 0.02350043838650813
$ws.Range("Q7").Value = 10.77171900582744
$ws.Range("R7").Value = 96.945471052447
$ws.Range("S7").Value = 0.02198026905927185
$ws.Range("T7").Value = 0.02198026905927185

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Comp"
$ws.Range("C8").Value = "Itgb3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.303350666666667
$ws.Range("H8").Value = 3.910052
$ws.Range("I8").Value = 0.04897540923545694
$ws.Range("J8").Value = 0.04897540923545693
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.970048
$ws.Range("N8").Value = 26.910144
$ws.Range("O8").Value = 0.487108783009476
$ws.Range("P8").Value = 0.4871087830094759
$ws.Range("Q8").Value = 11.691118040832
$ws.Range("R8").Value = 105.220062367488
$ws.Range("S8").Value = 0.02385635199007448
$ws.Range("T8").Value = 0.02385635199007447

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Comp"
$ws.Range("C9").Value = "Itgb3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.303350666666667
$ws.Range("H9").Value = 3.910052
$ws.Range("I9").Value = 0.04897540923545694
$ws.Range("J9").Value = 0.04897540923545693
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 9.012070666666666
$ws.Range("N9").Value = 27.036212
$ws.Range("O9").Value = 0.489390778604016
$ws.Range("P9").Value = 0.489390778604016
$ws.Range("Q9").Value = 11.74588831144711
$ws.Range("R9").Value = 105.712994803024
$ws.Range("S9").Value = 0.02396811365819058
$ws.Range("T9").Value = 0.02396811365819058

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Comp"
$ws.Range("C10").Value = "Itgb3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.303350666666667
$ws.Range("H10").Value = 3.910052
$ws.Range("I10").Value = 0.04897540923545694
$ws.Range("J10").Value = 0.04897540923545693
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4327576666666667
$ws.Range("N10").Value = 1.298273
$ws.Range("O10").Value = 0.02350043838650813
$ws.Range("P10").Value = 0.02350043838650813
$ws.Range("Q10").Value = 0.5640349933551111
$ws.Range("R10").Value = 5.076314940196
$ws.Range("S10").Value = 0.001150943587191877
$ws.Range("T10").Value = 0.001150943587191877
